# Update "gh-pages" generated data (广州-漫展信息.xlsx) — refresh counts scraped at 456a3b4.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 7817
$ws1.Range("F3").Value = 105
$ws1.Range("F4").Value = 83

# Row 5 (COMICUP 2024SP) went on sale: 想去人数 updated and 最低票价 changed
# from the placeholder text "不可售" to an actual numeric price; the cover
# image was also refreshed.
$ws1.Range("F5").Value = 11142
$ws1.Range("G5").Value = 68
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202403/niNTHKNX1711445933004.png"

$ws1.Range("F6").Value = 44
$ws1.Range("F8").Value = 645
$ws1.Range("F9").Value = 449
$ws1.Range("F12").Value = 778
$ws1.Range("F13").Value = 43
$ws1.Range("F14").Value = 83
$ws1.Range("F15").Value = 330
$ws1.Range("F16").Value = 19
$ws1.Range("F17").Value = 269
$ws1.Range("F19").Value = 397
$ws1.Range("F21").Value = 1091
$ws1.Range("F22").Value = 78
$ws1.Range("F23").Value = 630
$ws1.Range("F24").Value = 2223
$ws1.Range("F25").Value = 745
$ws1.Range("F27").Value = 554
$ws1.Range("F29").Value = 618
$ws1.Range("F30").Value = 554

# ---------------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 291
$ws2.Range("F3").Value = 66
$ws2.Range("F5").Value = 314
$ws2.Range("F11").Value = 3

# ---------------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 461

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) — union of the three sheets above
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 461
$ws4.Range("F3").Value = 7817
$ws4.Range("F4").Value = 105
$ws4.Range("F5").Value = 83
$ws4.Range("F6").Value = 291

# Row 7 mirrors 展览!F5/G5/I5 above (same event, same refresh).
$ws4.Range("F7").Value = 11154
$ws4.Range("G7").Value = 68
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202403/niNTHKNX1711445933004.png"

$ws4.Range("F8").Value = 44
$ws4.Range("F10").Value = 645
$ws4.Range("F11").Value = 449
$ws4.Range("F12").Value = 66
$ws4.Range("F16").Value = 314
$ws4.Range("F18").Value = 778
$ws4.Range("F19").Value = 43
$ws4.Range("F20").Value = 83
$ws4.Range("F21").Value = 330
$ws4.Range("F23").Value = 19
$ws4.Range("F27").Value = 269
$ws4.Range("F29").Value = 397
$ws4.Range("F31").Value = 1091
$ws4.Range("F32").Value = 78
$ws4.Range("F33").Value = 630
$ws4.Range("F34").Value = 2223
$ws4.Range("F35").Value = 745
$ws4.Range("F37").Value = 554
$ws4.Range("F39").Value = 3
$ws4.Range("F40").Value = 618
$ws4.Range("F41").Value = 554
